# Generate Report for Handoff
#
# Refreshes the latest-handoff timestamps for the "workerroletest" (zh-cn)
# and "Pong" (de-de) handoff batches, and fixes the Pong row's "Latest
# Handoff File" link on the zh-cn / de-de detail sheets, which had been
# erroneously pointing at the Ping handoff package instead of its own.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest Handoff Date" for workerroletest.md (row 4) and
# Pong.md (row 13) both refresh to the new handoff run timestamp.
$overview.Range("D4").Value  = "2016-04-04 03:46:26"
$overview.Range("D13").Value = "2016-04-04 03:46:26"

# zh-cn sheet, Pong.md row (row 13): the handoff file was mistakenly the
# Ping package - point it at the correct Pong package, and refresh the
# handoff datetime.
$zhcn.Range("D13").Value = "Pong.f5965988772320608f28a831c2d662c88665d3e1.zh-cn.xlf"
$zhcn.Range("E13").Value = "2016-04-04 03:46:19"

# zh-cn sheet, workerroletest.md row (row 4): refresh handoff datetime.
$zhcn.Range("E4").Value = "2016-04-04 03:46:19"

# de-de sheet, Pong.md row (row 11): same handoff-file fix and datetime
# refresh as zh-cn.
$dede.Range("D11").Value = "Pong.f5965988772320608f28a831c2d662c88665d3e1.de-de.xlf"
$dede.Range("E11").Value = "2016-04-04 03:46:26"

# de-de sheet, workerroletest.md row (row 21): refresh handback datetime.
$dede.Range("E21").Value = "2016-04-04 03:46:26"
